$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filtered save games) for rows 2-7, columns B-E and G.
# Column F (Win) is unchanged.

$data = @{
    2 = @{ B = 0.1190320826869504;    C = 0.306821227259698;    D = 3.537761648806719;  E = 0.4942365360607697; G = 4.457851494814137 }
    3 = @{ B = 0.0006408296065709695; C = 0.04071648406533734;  D = 6708.013860684405;  E = 2195978.878461985;  G = 2202686.933679983 }
    4 = @{ B = 1.455362044514542;     C = 1.655778082260271;    D = 0.7527432677738641; E = 0.4942365360607697; G = 4.358119930609447 }
    5 = @{ B = 3.286832544864788;     C = 10.34677158129881;    D = 3.537761648806719;  E = 10.19245300693656;  G = 27.36381878190688 }
    6 = @{ B = 1.455362044514542;     C = 1.655778082260271;    D = 0.7527432677738641; E = 0.4942365360607697; G = 4.358119930609447 }
    7 = @{ B = 0.2917716402565462;    C = 0.306821227259698;    D = 0.1494219747398047; E = 0.4942365360607697; G = 1.242251378316819 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
